$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G4').Value = 'This indicator measures the level of engagement and support provided by Rotarians in WASH-related events in local communes, reflecting increased civil society engagement in WASH services.'
$ws.Range('I4').Value = 38
$ws.Range('K4').Value = 65

$ws.Range('G5').Value = 'The percentage of WASH events within communes that are led and coordinated by trained HANWASH Ambassadors and Coordinators, aiming to effectively communicate WASH messages and advocate for HANWASH goals. This reflects the level of leadership and initiative taken by trained personnel in promoting WASH activities.'
$ws.Range('I5').Value = 71
$ws.Range('K5').Value = 94

$ws.Range('G6').Value = 'The total number of hours that Rotarians volunteer each month participating in WASH initiatives and activities.'
$ws.Range('I6').Value = 9
$ws.Range('J6').Value = 39
$ws.Range('K6').Value = 27

$ws.Range('G7').Value = 'The total number of individuals who participate in or are impacted by WASH events and advocacy efforts led by HANWASH Ambassadors, Coordinators, and Rotarians. This indicates the outreach and engagement success of the WASH programs.'
$ws.Range('I7').Value = 15
$ws.Range('J7').Value = 38
$ws.Range('K7').Value = 21

$ws.Range('G8').Value = 'The total number of HANWASH Ambassadors who have successfully completed their training programs.'
$ws.Range('I8').Value = 5
$ws.Range('J8').Value = 75
$ws.Range('K8').Value = 60

$ws.Range('G9').Value = 'The total number of HANWASH Coordinators who have successfully completed their training programs.'
$ws.Range('I9').Value = 17
$ws.Range('J9').Value = 77
$ws.Range('K9').Value = 37

$ws.Range('G10').Value = ' The total count of educational and training materials developed that adhere to the three pillars and core values of HANWASH.'
$ws.Range('I10').Value = 11
$ws.Range('J10').Value = 12
$ws.Range('K10').Value = 12

$ws.Range('I11').Value = 23
$ws.Range('J11').Value = 34
$ws.Range('K11').Value = 25

$ws.Range('G12').Value = ' The percentage of WASH interventions that have been approved and implemented in accordance with predefined Commune Action Plans, ensuring the initiatives are part of a strategic and coordinated approach to WASH improvements.'
$ws.Range('I12').Value = 67
$ws.Range('K12').Value = 98

$ws.Range('G13').Value = 'The number of communes that have developed action plans explicitly addressing water, sanitation, hygiene, and water resources management.'
$ws.Range('I13').Value = 1
$ws.Range('J13').Value = 1
$ws.Range('K13').Value = 1

$ws.Range('G14').Value = 'The number of meetings and events held to plan and coordinate WASH activities involving DINEPA and municipal officials.'
$ws.Range('J14').Value = 1
$ws.Range('K14').Value = 1

$ws.Range('G15').Value = ' The percentage of service providers participating in the HANWASH initiative that are regularly monitored according to established national guidelines. This metric ensures that WASH service providers meet accountability standards, maintain quality service delivery, and uphold their responsibilities towards users as per national regulatory requirements.'
$ws.Range('I15').Value = 33
$ws.Range('K15').Value = 37

$ws.Range('G16').Value = 'The percentage of service providers with established accountability mechanisms, such as regular meetings, community engagement, bylaws, fee collection, record-keeping, and monthly reporting to local authorities.'
$ws.Range('I16').Value = 67
$ws.Range('K16').Value = 87

$ws.Range('I17').Value = 32
$ws.Range('K17').Value = 91

$ws.Range('G18').Value = 'The number of water systems within intervention areas that consistently submit monthly operational and performance reports, ensuring regular monitoring and transparency of water service delivery.'

$ws.Range('G19').Value = ' The percentage of customer complaints that are resolved within the same month they are identified.'
$ws.Range('I19').Value = 74
$ws.Range('K19').Value = 97

$ws.Range('I20').Value = 7
$ws.Range('J20').Value = 11
$ws.Range('K20').Value = 8

$ws.Range('G22').Value = 'The percentage of water points in intervention areas that are operational (either fully functional or functional with minor repairs needed) and provide potable water, free from E. Coli and other priority contaminants, according to original design specifications.'
$ws.Range('I22').Value = 46
$ws.Range('K22').Value = 97

$ws.Range('G23').Value = 'The percentage of people in intervention communes who have access to an improved water source within a 30-minute round trip, ensuring basic drinking water service at the household level.'
$ws.Range('I23').Value = 71
$ws.Range('K23').Value = 74

$ws.Range('G24').Value = 'The percentage of people in target areas who have access to safe drinking water that is always available, accessible at home, and free from contaminants.'
$ws.Range('I24').Value = 78
$ws.Range('K24').Value = 81

$ws.Range('G25').Value = 'The percentage of water points in intervention areas that are fully operational, provide potable water, and maintain a balanced or surplus budget within two years from their inauguration date.'
$ws.Range('I25').Value = 47
$ws.Range('K25').Value = 93

$ws.Range('G26').Value = 'The number of individuals in intervention communes who benefit from access to basic drinking water services provided by community-managed water points, ensuring sustainable and locally-managed water supply solutions.'
$ws.Range('I26').Value = 51
$ws.Range('J26').Value = 95
$ws.Range('K26').Value = 63

$ws.Range('G27').Value = 'The total number of water points established and managed by local communities under DINEPA''s guidelines.'
$ws.Range('I27').Value = 35
$ws.Range('J27').Value = 185
$ws.Range('K27').Value = 72

$ws.Range('G28').Value = 'The percentage of community-managed water points that submit monthly reports to local authorities.'
$ws.Range('I28').Value = 70
$ws.Range('J28').Value = 164
$ws.Range('K28').Value = 113

$ws.Range('G29').Value = 'The number of individuals who have gained access to safe drinking water through the project.'
$ws.Range('I29').Value = 40
$ws.Range('J29').Value = 92

$ws.Range('G30').Value = 'The number of existing piped water systems that have been restored to proper functionality in the target areas.'
$ws.Range('J30').Value = 1
$ws.Range('K30').Value = 1

$ws.Range('G31').Value = 'The number of new piped water systems established and professionally managed in the target areas.'

$ws.Range('G32').Value = 'The percentage of communities in intervention areas that have been verified as Open Defecation Free by the commune WASH committee, meeting criteria such as the presence of usable toilets, no fecal matter in open areas, and willingness to maintain ODF status.'
$ws.Range('I32').Value = 89
$ws.Range('K32').Value = 97

$ws.Range('G33').Value = ' The percentage of people in intervention communes who have access to improved and unshared sanitation facilities, ensuring basic sanitation service at the household level.'
$ws.Range('I33').Value = 73
$ws.Range('K33').Value = 83

$ws.Range('G34').Value = 'The percentage of people in target areas who have access to safe, private (unshared) sanitation facilities where excreta are properly treated on-site or transported and treated off-site.'
$ws.Range('I34').Value = 51
$ws.Range('K34').Value = 81

$ws.Range('G35').Value = 'The number of individuals newly provided with access to basic sanitation services, including improved and unshared facilities, within intervention communes.'
$ws.Range('I35').Value = 9
$ws.Range('J35').Value = 12
$ws.Range('K35').Value = 11

$ws.Range('G36').Value = 'The total number of communities in intervention areas that have achieved verification as Open Defecation Free, meeting the criteria set by the commune WASH committee.'
$ws.Range('I36').Value = 21
$ws.Range('J36').Value = 140
$ws.Range('K36').Value = 125

$ws.Range('G37').Value = 'The number of communities in intervention areas that have sustained Open Defecation Free status for at least one year and have been officially certified by the commune WASH committee.'
$ws.Range('I37').Value = 86
$ws.Range('J37').Value = 133
$ws.Range('K37').Value = 103

$ws.Range('G38').Value = 'The number of public latrines built within the project scope.'
$ws.Range('I38').Value = 9
$ws.Range('J38').Value = 106
$ws.Range('K38').Value = 43

$ws.Range('G39').Value = 'The percentage of people in intervention communes who have access to facilities with soap and water for handwashing, ensuring basic hygiene service at the household level.'
$ws.Range('I39').Value = 20
$ws.Range('K39').Value = 95

$ws.Range('G40').Value = 'The percentage of individuals newly provided with access to basic hygiene services, including handwashing facilities with soap and water within intervention communes.'
$ws.Range('I40').Value = 98
$ws.Range('K40').Value = 99

$ws.Range('G41').Value = 'The total number of community animators who have completed training programs.'
$ws.Range('I41').Value = 85
$ws.Range('J41').Value = 95
$ws.Range('K41').Value = 90

$ws.Range('I42').Value = 14
$ws.Range('J42').Value = 19
$ws.Range('K42').Value = 17

$ws.Range('G43').Value = 'The number of individuals newly provided with access to basic hygiene services, including handwashing facilities with soap and water within intervention communes.'
$ws.Range('I43').Value = 41
$ws.Range('J43').Value = 46
$ws.Range('K43').Value = 42

$ws.Range('G44').Value = 'The percentage of schools in intervention areas that have improved, safe drinking water sources, functional sanitation facilities, and handwashing stations with soap and water.'
$ws.Range('I44').Value = 45
$ws.Range('K44').Value = 84

$ws.Range('G45').Value = 'The percentage of healthcare facilities in intervention areas that have improved, safe drinking water sources, functioning sanitation facilities, and handwashing stations with soap and water or alcohol-based hand rubs. This metric ensures that healthcare facilities meet basic WASH standards.'
$ws.Range('I45').Value = 81
$ws.Range('K45').Value = 88

$ws.Range('G46').Value = 'The total number of schools in intervention areas that have been provided with access to improved drinking water services, ensuring that they meet the required standards of water quality and availability.'
$ws.Range('I46').Value = 61
$ws.Range('J46').Value = 168
$ws.Range('K46').Value = 148

$ws.Range('G47').Value = 'The total number of healthcare facilities in intervention areas that have been provided with access to improved drinking water services, ensuring safe and reliable water for patients and staff.'
$ws.Range('I47').Value = 175
$ws.Range('J47').Value = 180
$ws.Range('K47').Value = 176

$ws.Range('G48').Value = 'The total number of schools in intervention areas that have been equipped with improved, functional sanitation facilities, including gender-separated and accessible toilets, to meet the needs of students and staff.'
$ws.Range('I48').Value = 116
$ws.Range('J48').Value = 183
$ws.Range('K48').Value = 131

$ws.Range('G49').Value = 'The total number of healthcare facilities in intervention areas that have been equipped with improved, functional sanitation facilities, including dedicated and accessible toilets. These facilities are designed to meet the sanitation needs of all healthcare users.'
$ws.Range('I49').Value = 106
$ws.Range('J49').Value = 162
$ws.Range('K49').Value = 133

$ws.Range('G50').Value = 'The total number of schools in intervention areas that have been provided with handwashing facilities, including soap and water, ensuring basic hygiene practices are accessible to all students and staff.'
$ws.Range('I50').Value = 169
$ws.Range('J50').Value = 194
$ws.Range('K50').Value = 192

$ws.Range('G51').Value = 'The total number of healthcare facilities in intervention areas that have been provided with handwashing facilities, including soap and water or alcohol-based hand rubs, ensuring basic hygiene practices are accessible to all patients and staff.'
$ws.Range('I51').Value = 14
$ws.Range('J51').Value = 105
$ws.Range('K51').Value = 97

$ws.Range('G53').Value = 'The total amount of funds pledged by stakeholders in accordance with HANWASH Core Values.'
$ws.Range('I53').Value = 646
$ws.Range('J53').Value = 10675
$ws.Range('K53').Value = 2899

$ws.Range('G54').Value = 'The cumulative percentage of allocated funds that have been spent over the project duration. This metric is calculated by dividing the total amount spent by the total committed funds and multiplying by 100, reflecting the financial resource utilization efficiency of the project.'
$ws.Range('I54').Value = 96
$ws.Range('K54').Value = 98

$ws.Range('G55').Value = 'The total amount of money spent by external organizations or entities within the areas covered by HANWASH projects.'
$ws.Range('I55').Value = 5376
$ws.Range('J55').Value = 5669
$ws.Range('K55').Value = 5525

$ws.Range('G56').Value = 'The total amount of money spent by external entities outside the HANWASH project areas, in line with HANWASH Core Values.'
$ws.Range('I56').Value = 7859
$ws.Range('J56').Value = 15348
$ws.Range('K56').Value = 13865

$ws.Range('I57').Value = 74
$ws.Range('K57').Value = 77

$ws.Range('G58').Value = 'The total number of DINEPA personnel who have undergone leadership training programs.'
$ws.Range('I58').Value = 33
$ws.Range('J58').Value = 44
$ws.Range('K58').Value = 42

$ws.Range('G59').Value = 'The number of technical training sessions delivered to address DINEPA''s prioritized areas, including the unified national tariff methodology and ODF certification for 2023-24.'
$ws.Range('I59').Value = 24
$ws.Range('J59').Value = 83
$ws.Range('K59').Value = 68
